# Refatorando artefatos seguindo a correção da ac5
# Remove the last row of the restrictions table (the row describing the
# client's request for additional features and the corresponding answer),
# as that constraint/reason pair no longer applies.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Locate the row whose first cell starts with the text that identifies it,
# then delete the whole row (both cells included).
for ($i = $table.Rows.Count; $i -ge 1; $i--) {
    $row = $table.Rows.Item($i)
    $cellText = $row.Cells.Item(1).Range.Text
    if ($cellText -like "*Cliente solicitou que seja entregue mais funcionalidades*") {
        $row.Delete()
        break
    }
}
